# Update crypto price/volume data per the Mon Sep 25 16:16:59 UTC 2023 GitHub Actions run.
# Cells in columns D (Price) and E (Volume 1h), plus the B/C identity swap for
# rows 43-44 (Aave / MXToken traded ranking places), are refreshed in place.
#
# Set-CellText forces the destination cell to Text before assigning the
# literal string and then clears the temporary number-format override, so
# values that look numeric (e.g. "0.999", "209.85") are stored as plain text
# -- matching the workbook's existing convention for this sheet -- instead of
# being auto-converted to numbers/dates by Excel's normal type inference.
function Set-CellText($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws 'D2' '26.330.80'
Set-CellText $ws 'E2' '  -1.50%  '
Set-CellText $ws 'D3' '1.592.41'
Set-CellText $ws 'E3' '  -0.41%  '
Set-CellText $ws 'D4' '0.999'
Set-CellText $ws 'E4' '  -0.66%  '
Set-CellText $ws 'D5' '209.85'
Set-CellText $ws 'E5' '  -0.93%  '
Set-CellText $ws 'D6' '0.501'
Set-CellText $ws 'E6' '  -2.22%  '
Set-CellText $ws 'D7' '0.999'
Set-CellText $ws 'E7' '  -0.67%  '
Set-CellText $ws 'D8' '0.0612'
Set-CellText $ws 'E8' '  -1.06%  '
Set-CellText $ws 'D9' '0.246'
Set-CellText $ws 'E9' '  -0.42%  '
Set-CellText $ws 'D10' '19.66'
Set-CellText $ws 'E10' '  +0.04%  '
Set-CellText $ws 'D11' '0.0845'
Set-CellText $ws 'E11' '  -0.29%  '
Set-CellText $ws 'D12' '1.813.36'
Set-CellText $ws 'E12' '  -0.60%  '
Set-CellText $ws 'D13' '1.630.24'
Set-CellText $ws 'E13' '  +1.30%  '
Set-CellText $ws 'D14' '4.08'
Set-CellText $ws 'E14' '  +0.41%  '
Set-CellText $ws 'D15' '0.519'
Set-CellText $ws 'E15' '  -1.19%  '
Set-CellText $ws 'D16' '64.84'
Set-CellText $ws 'E16' '  -0.35%  '
Set-CellText $ws 'D17' '26.308.98'
Set-CellText $ws 'E17' '  -1.48%  '
Set-CellText $ws 'D18' '0.0₃0730'
Set-CellText $ws 'E18' '  -1.33%  '
Set-CellText $ws 'D19' '7.43'
Set-CellText $ws 'E19' '  +3.83%  '
Set-CellText $ws 'D20' '213.10'
Set-CellText $ws 'E20' '  +2.14%  '
Set-CellText $ws 'E21' '  -0.55%  '
Set-CellText $ws 'D22' '4.28'
Set-CellText $ws 'E22' '  -0.33%  '
Set-CellText $ws 'D23' '2.18'
Set-CellText $ws 'E23' '  -2.30%  '
Set-CellText $ws 'D24' '8.90'
Set-CellText $ws 'E24' '  -1.50%  '
Set-CellText $ws 'D25' '144.64'
Set-CellText $ws 'E25' '  +0.48%  '
Set-CellText $ws 'D27' '7.05'
Set-CellText $ws 'E27' '  -1.23%  '
Set-CellText $ws 'E28' '  -1.22%  '
Set-CellText $ws 'D29' '15.32'
Set-CellText $ws 'E29' '  -0.28%  '
Set-CellText $ws 'E30' '  -0.16%  '
Set-CellText $ws 'E31' '  -0.94%  '
Set-CellText $ws 'E32' '  -1.18%  '
Set-CellText $ws 'E33' '  +1.06%  '
Set-CellText $ws 'D34' '1.291.59'
Set-CellText $ws 'E34' '  +1.31%  '
Set-CellText $ws 'D35' '2.44'
Set-CellText $ws 'E35' '  -1.71%  '
Set-CellText $ws 'D36' '0.608'
Set-CellText $ws 'E36' '  +2.81%  '
Set-CellText $ws 'D37' '1.48'
Set-CellText $ws 'E37' '  -0.96%  '
Set-CellText $ws 'D38' '1.12'
Set-CellText $ws 'E38' '  -9.04%  '
Set-CellText $ws 'E39' '  -1.29%  '
Set-CellText $ws 'D40' '0.814'
Set-CellText $ws 'E40' '  -1.34%  '
Set-CellText $ws 'D41' '0.999'
Set-CellText $ws 'E41' '  -0.69%  '
Set-CellText $ws 'E42' '  +3.09%  '
Set-CellText $ws 'B43' 'Aave'
Set-CellText $ws 'C43' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-CellText $ws 'D43' '62.79'
Set-CellText $ws 'E43' '  +0.22%  '
Set-CellText $ws 'B44' 'MXToken'
Set-CellText $ws 'C44' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText $ws 'D44' '2.14'
Set-CellText $ws 'E44' '  -2.45%  '
Set-CellText $ws 'E45' '  -1.99%  '
Set-CellText $ws 'D46' '1.726.55'
Set-CellText $ws 'E46' '  -0.54%  '
Set-CellText $ws 'D47' '88.71'
Set-CellText $ws 'E47' '  -1.95%  '
Set-CellText $ws 'D48' '1.51'
Set-CellText $ws 'E48' '  -3.57%  '
Set-CellText $ws 'D49' '0.0996'
Set-CellText $ws 'E49' '  -2.65%  '
Set-CellText $ws 'E50' '  -1.35%  '
Set-CellText $ws 'D51' '0.0₇0981'
Set-CellText $ws 'E51' '  -7.11%  '
